$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'capri hex',
    'knee pads softball',
    'wrestling knee pads mens',
    'snowboarding knee brace',
    'basketball leggings for boys youth',
    'honeycomb leggings',
    'sliding knee pad',
    'compression knee pad pants',
    'youth basketball compression pants',
    'basketball guide hand',
    '28 inch basketball',
    'mens tights with knee pads',
    'sliding sleeve softball',
    'knee sleeve youth wrestling',
    'basketball pads knee',
    'sliding pants baseball youth',
    'sports knee pads basketball',
    'knee pads crossfit',
    'boys knee pad leggings',
    'knee pad basketball',
    'basketball pads youth',
    'knee pads for snowboarding',
    'basketball leg sleeves for boys',
    'basketball leggings youth boys',
    '20 30 compression leggings',
    'workout knee pads',
    'honeycomb compression pants',
    'boys basketball leg sleeve',
    'knee brace hockey',
    'youth thigh compression sleeve',
    'basketball tights youth boys',
    'boys volleyball knee pads',
    'crossfit pads',
    'compression knee pads youth',
    'padded basketball leggings for boys',
    'crossfit knee sleeves men',
    'baseball pants mens knee high',
    'basketball knee pads for youth',
    'softball knee pad',
    'kneepads men',
    'knee compression sleeve reduce strain',
    'knee brace',
    'knee compression sleeve',
    'knee support',
    'knee brace support',
    'knee brace compression sleeve',
    'compression sleeve',
    'knee sleeve',
    'knee compression',
    'compression knee brace',
    'knee pad',
    'thigh compression sleeve',
    'basketball training',
    'crossfit training',
    'basketball thigh pads',
    'waist training leggings',
    'work out tights mens',
    'compression leggings capri',
    'black workout leggings',
    'workout legging',
    'capri legging',
    'workout tights men',
    'workout pad',
    'leggings basketball',
    'black legging',
    'basketball knee sleeve men',
    'basketball pads for training',
    'workout tights for men',
    'knee pads basketball mcdavid black',
    'basketball knee pads adult black',
    'mens leggings compression',
    'black capri legging',
    'man leggings',
    'basketball sleeve black',
    'men workout pants',
    'basketball knee pad',
    'hip workout',
    'compression legging',
    'workout pant',
    'capri pant',
    'youth training basketball',
    'compression pant',
    'baseball pants knee high mens',
    'basketball knee pads youth girls',
    'basketball knee protector',
    'basketball leggings with knee pads',
    'basketball youth knee pads',
    'compression knee pants',
    'compression leggings boys basketball',
    'knee guards for basketball',
    'knee pad basketball men',
    'knee pad pants basketball',
    'knee pads basketball',
    'knee pads boys basketball',
    'knee pads compression pants',
    'knee pads for basketball men',
    'knee pads men basketball',
    'knee pads pants for men',
    'knee pads wrestling',
    'knee protection for workout'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

